# Auto-generated edit script: updates crypto price/volume data
# per the Mon Oct  2 06:51:42 UTC 2023 GitHub Actions refresh commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.064.94'
$ws.Range('E2').Value = '  +3.43%  '
$ws.Range('D3').Value = '1.723.17'
$ws.Range('E3').Value = '  +2.59%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '218.93'
$ws.Range('E5').Value = '  +1.82%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '24.54'
$ws.Range('E8').Value = '  +14.78%  '
$ws.Range('E9').Value = '  +3.36%  '
$ws.Range('E11').Value = '  +1.36%  '
$ws.Range('D12').Value = '1.967.06'
$ws.Range('E12').Value = '  +2.67%  '
$ws.Range('D13').Value = '1.723.22'
$ws.Range('E13').Value = '  +2.50%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.27'
$ws.Range('E14').Value = '  +2.99%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.563'
$ws.Range('E15').Value = '  +4.95%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '67.59'
$ws.Range('E16').Value = '  +2.07%  '
$ws.Range('D17').Value = '28.029.51'
$ws.Range('E17').Value = '  +3.38%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '243.11'
$ws.Range('E18').Value = '  +1.91%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '8.03'
$ws.Range('E19').Value = '  -1.32%  '
$ws.Range('E20').Value = '  +1.41%  '
$ws.Range('E21').Value = '  -0.03%  '
$ws.Range('E22').Value = '  +2.48%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.72'
$ws.Range('E23').Value = '  +2.72%  '
$ws.Range('E24').Value = '  -0.13%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '149.04'
$ws.Range('E25').Value = '  +1.58%  '
$ws.Range('E26').Value = '  +3.42%  '
$ws.Range('E27').Value = '  +2.53%  '
$ws.Range('E28').Value = '  +0.94%  '
$ws.Range('E29').Value = '  +0.12%  '
$ws.Range('E30').Value = '  +2.07%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.20'
$ws.Range('E31').Value = '  +2.00%  '
$ws.Range('E32').Value = '  +2.11%  '
$ws.Range('D33').Value = '1.495.89'
$ws.Range('E33').Value = '  -3.97%  '
$ws.Range('E34').Value = '  +2.17%  '
$ws.Range('E35').Value = '  -2.28%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.957'
$ws.Range('E36').Value = '  +2.28%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.607'
$ws.Range('E37').Value = '  +0.86%  '
$ws.Range('E38').Value = '  +0.59%  '
$ws.Range('E39').Value = '  +0.32%  '
$ws.Range('E40').Value = '  +0.96%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '70.77'
$ws.Range('E41').Value = '  +2.01%  '
$ws.Range('E43').Value = '  -0.02%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.28'
$ws.Range('E44').Value = '  +1.50%  '
$ws.Range('D45').Value = '1.870.86'
$ws.Range('E45').Value = '  +2.50%  '
$ws.Range('E46').Value = '  +2.74%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.77'
$ws.Range('E47').Value = '  +11.57%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '90.98'
$ws.Range('E48').Value = '  +0.27%  '
$ws.Range('D49').Value = '0.0₆0113'
$ws.Range('E49').Value = '  +5.77%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.26'
$ws.Range('E50').Value = '  +1.57%  '
$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.105'
$ws.Range('E51').Value = '  +0.62%  '

Write-Output "Applied 77 cell updates"
